# Apply updated betting-odds values (Jogos da Semana FlashScore 2025-02-11)
# Cells are set per data row; values below match the published diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.25
$ws.Range("H2").Value = 2.6
$ws.Range("I2").Value = 2.75
$ws.Range("L2").Value = 3.6
$ws.Range("N2").Value = 4.75
$ws.Range("AC2").Value = 6.5
$ws.Range("AI2").Value = 4.75

# Row 3
$ws.Range("H3").Value = 3.6
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 2.35
$ws.Range("T3").Value = 1.57
$ws.Range("U3").Value = 3.85
$ws.Range("V3").Value = 1.27

# Row 6
$ws.Range("G6").Value = 2.9
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 3.6
$ws.Range("L6").Value = 3.5
$ws.Range("AE6").Value = 11
$ws.Range("AG6").Value = 26
$ws.Range("AO6").Value = 12
$ws.Range("AQ6").Value = 29
$ws.Range("AR6").Value = 26

# Row 7
$ws.Range("G7").Value = 1.18
$ws.Range("H7").Value = 6.3
$ws.Range("I7").Value = 12
$ws.Range("J7").Value = 1.53
$ws.Range("K7").Value = 2.77
$ws.Range("L7").Value = 9.25
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.47
$ws.Range("W7").Value = 2
$ws.Range("X7").Value = 1.65
$ws.Range("AB7").Value = 1.6
$ws.Range("AC7").Value = 8.75
$ws.Range("AD7").Value = 6.4
$ws.Range("AE7").Value = 10
$ws.Range("AF7").Value = 6.8
$ws.Range("AH7").Value = 32
$ws.Range("AI7").Value = 17.5
$ws.Range("AJ7").Value = 13.5
$ws.Range("AK7").Value = 30
$ws.Range("AL7").Value = 150
$ws.Range("AN7").Value = 35
$ws.Range("AO7").Value = 100
$ws.Range("AP7").Value = 40
$ws.Range("AQ7").Value = 450
$ws.Range("AR7").Value = 175

# Row 8
$ws.Range("AD8").Value = 15
$ws.Range("AE8").Value = 12
$ws.Range("AF8").Value = 34

# Row 10
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 2.6
$ws.Range("J10").Value = 3.7
$ws.Range("K10").Value = 1.8
$ws.Range("L10").Value = 3.55
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 5.1
$ws.Range("O10").Value = 1.52
$ws.Range("P10").Value = 2.35
$ws.Range("S10").Value = 2.55
$ws.Range("T10").Value = 1.45
$ws.Range("W10").Value = 4.5
$ws.Range("X10").Value = 1.17
$ws.Range("Z10").Value = 2.2
$ws.Range("AA10").Value = 1.98
$ws.Range("AC10").Value = 7.1
$ws.Range("AD10").Value = 14.5
$ws.Range("AE10").Value = 10.75
$ws.Range("AF10").Value = 40
$ws.Range("AG10").Value = 30
$ws.Range("AI10").Value = 5.1
$ws.Range("AJ10").Value = 5.2
$ws.Range("AN10").Value = 6.5
$ws.Range("AO10").Value = 12.5
$ws.Range("AR10").Value = 29

# Row 11
$ws.Range("G11").Value = 3.55
$ws.Range("H11").Value = 2.87
$ws.Range("I11").Value = 2.2
$ws.Range("J11").Value = 4.05
$ws.Range("K11").Value = 1.95
$ws.Range("L11").Value = 2.87
$ws.Range("N11").Value = 5.7
$ws.Range("P11").Value = 2.57
$ws.Range("S11").Value = 2.3
$ws.Range("T11").Value = 1.55
$ws.Range("W11").Value = 4
$ws.Range("X11").Value = 1.2
$ws.Range("Y11").Value = 1.5
$ws.Range("Z11").Value = 2.42
$ws.Range("AC11").Value = 8.75
$ws.Range("AD11").Value = 18.5
$ws.Range("AE11").Value = 12
$ws.Range("AF11").Value = 55
$ws.Range("AG11").Value = 35
$ws.Range("AH11").Value = 45
$ws.Range("AI11").Value = 5.7
$ws.Range("AJ11").Value = 5.6
$ws.Range("AK11").Value = 15
$ws.Range("AO11").Value = 9.5
$ws.Range("AP11").Value = 9
$ws.Range("AQ11").Value = 22
$ws.Range("AR11").Value = 21
$ws.Range("AS11").Value = 37

# Row 12
$ws.Range("G12").Value = 1.26
$ws.Range("H12").Value = 4.9
$ws.Range("J12").Value = 1.72
$ws.Range("K12").Value = 2.35
$ws.Range("N12").Value = 7.9
$ws.Range("O12").Value = 1.25
$ws.Range("P12").Value = 3.6
$ws.Range("S12").Value = 1.75
$ws.Range("T12").Value = 1.98
$ws.Range("W12").Value = 2.77
$ws.Range("X12").Value = 1.39
$ws.Range("Y12").Value = 1.38
$ws.Range("Z12").Value = 2.82
$ws.Range("AA12").Value = 2.35
$ws.Range("AB12").Value = 1.53
$ws.Range("AC12").Value = 6
$ws.Range("AD12").Value = 5.4
$ws.Range("AG12").Value = 11.75
$ws.Range("AI12").Value = 7.9
$ws.Range("AJ12").Value = 10.25
$ws.Range("AN12").Value = 27

# Row 13
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.75
$ws.Range("S13").Value = 1.98
$ws.Range("T13").Value = 1.88

# Row 14
$ws.Range("W14").Value = 5
$ws.Range("X14").Value = 1.17

# Row 15
$ws.Range("G15").Value = 2.15
$ws.Range("K15").Value = 2.1
$ws.Range("L15").Value = 3.75
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.3
$ws.Range("P15").Value = 3.5
$ws.Range("S15").Value = 2.03
$ws.Range("T15").Value = 1.83
$ws.Range("W15").Value = 3.5
$ws.Range("X15").Value = 1.3
$ws.Range("Y15").Value = 1.4
$ws.Range("Z15").Value = 2.75
$ws.Range("AA15").Value = 1.73
$ws.Range("AB15").Value = 2
$ws.Range("AC15").Value = 8
$ws.Range("AD15").Value = 11
$ws.Range("AE15").Value = 9
$ws.Range("AG15").Value = 17
$ws.Range("AH15").Value = 26
$ws.Range("AI15").Value = 10
$ws.Range("AK15").Value = 13
$ws.Range("AL15").Value = 41
$ws.Range("AM15").Value = 201
$ws.Range("AN15").Value = 10
$ws.Range("AO15").Value = 17
$ws.Range("AR15").Value = 26
$ws.Range("AS15").Value = 34

# Row 16
$ws.Range("M16").Value = 1.05
$ws.Range("N16").Value = 11
$ws.Range("O16").Value = 1.25
$ws.Range("P16").Value = 4
$ws.Range("S16").Value = 1.8
$ws.Range("T16").Value = 2
$ws.Range("W16").Value = 3
$ws.Range("X16").Value = 1.4

# Row 18
$ws.Range("M18").Value = 1.03
$ws.Range("N18").Value = 15

# Row 19
$ws.Range("G19").Value = 2.1
$ws.Range("J19").Value = 2.88
$ws.Range("AD19").Value = 10
$ws.Range("AG19").Value = 19
$ws.Range("AS19").Value = 34

# Row 20
$ws.Range("G20").Value = 5.75
$ws.Range("H20").Value = 3.7
$ws.Range("I20").Value = 1.57
$ws.Range("N20").Value = 9.5
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 1.85
$ws.Range("AI20").Value = 9.5
$ws.Range("AO20").Value = 7

# Row 21
$ws.Range("G21").Value = 1.8
$ws.Range("H21").Value = 3.9
$ws.Range("I21").Value = 3.8
$ws.Range("J21").Value = 2.4
$ws.Range("K21").Value = 2.25
$ws.Range("L21").Value = 4.5
$ws.Range("M21").Value = 1.04
$ws.Range("N21").Value = 13
$ws.Range("O21").Value = 1.25
$ws.Range("P21").Value = 3.75
$ws.Range("S21").Value = 1.75
$ws.Range("T21").Value = 2.05
$ws.Range("W21").Value = 3
$ws.Range("X21").Value = 1.36
$ws.Range("Y21").Value = 1.36
$ws.Range("Z21").Value = 3
$ws.Range("AD21").Value = 9
$ws.Range("AF21").Value = 15
$ws.Range("AG21").Value = 15
$ws.Range("AI21").Value = 13
$ws.Range("AJ21").Value = 7.5
$ws.Range("AN21").Value = 12
$ws.Range("AO21").Value = 21
$ws.Range("AP21").Value = 13
$ws.Range("AQ21").Value = 41
$ws.Range("AR21").Value = 29

# Row 23
$ws.Range("G23").Value = 1.83
$ws.Range("H23").Value = 3.6
$ws.Range("I23").Value = 4
$ws.Range("AG23").Value = 17
$ws.Range("AI23").Value = 9
$ws.Range("AM23").Value = 401
$ws.Range("AP23").Value = 15

# Row 24
$ws.Range("I24").Value = 3.3
$ws.Range("M24").Value = 1.06
$ws.Range("N24").Value = 10
$ws.Range("AN24").Value = 9.5
$ws.Range("AS24").Value = 34

# Row 26
$ws.Range("H26").Value = 3.3
$ws.Range("K26").Value = 2.1
$ws.Range("M26").Value = 1.05
$ws.Range("N26").Value = 11
$ws.Range("S26").Value = 1.98
$ws.Range("T26").Value = 1.83

# Row 27
$ws.Range("G27").Value = 8.25
$ws.Range("I27").Value = 1.27
$ws.Range("J27").Value = 7.1
$ws.Range("K27").Value = 2.65
$ws.Range("O27").Value = 1.16
$ws.Range("P27").Value = 4.65
$ws.Range("S27").Value = 1.5
$ws.Range("T27").Value = 2.42
$ws.Range("W27").Value = 2.18
$ws.Range("X27").Value = 1.62
$ws.Range("AA27").Value = 1.88
$ws.Range("AB27").Value = 1.82
$ws.Range("AC27").Value = 26
$ws.Range("AD27").Value = 65
$ws.Range("AE27").Value = 27
$ws.Range("AF27").Value = 300
$ws.Range("AJ27").Value = 11
$ws.Range("AM27").Value = 700
$ws.Range("AN27").Value = 8.5
$ws.Range("AO27").Value = 6.9
$ws.Range("AP27").Value = 8.75
$ws.Range("AS27").Value = 25
